# Add a new "localdb" command-category column to the hidden '#system' sheet,
# used by the "filter" sheet's dropdown (INDIRECT(target-category)) to list
# nexial command signatures per category.
#
# Layout recap (row 1 = category name, rows 2..N = command signatures):
#   Column A "target"  -> alphabetical list of every category name
#   Columns B..AC       -> one category per column (alphabetical by header)
#
# "localdb" sorts alphabetically right after "json" and before "macro", i.e.
# at column N (14) -- so every existing column from N onward shifts right by
# one letter, and "localdb" is inserted into the "target" list (column A)
# between "json" and "macro" as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) Make room: insert a blank column before N (14). Everything N..AC moves
#    to O..AD, formulas/defined names aside (handled explicitly below).
$ws.Columns.Item(14).Insert()

# 2) "target" (column A) keeps every category name in alphabetical order.
#    Shift existing entries for rows 14..29 ("macro" onward) down one row,
#    bottom-up so we don't clobber data before reading it.
for ($r = 29; $r -ge 14; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value()
}
$ws.Cells.Item(14, 1).Value = "localdb"

# 3) Populate the new "localdb" column (N) header + its six commands.
$ws.Cells.Item(1, 14).Value = "localdb"
$ws.Cells.Item(2, 14).Value = "cloneTable(var,source,target)"
$ws.Cells.Item(3, 14).Value = "dropTables(var,tables)"
$ws.Cells.Item(4, 14).Value = "exportCSV(sql,output)"
$ws.Cells.Item(5, 14).Value = "importRecords(var,sourceDb,sql,table)"
$ws.Cells.Item(6, 14).Value = "purge(var)"
$ws.Cells.Item(7, 14).Value = "runSQLs(var,sqls)"

# 4) Re-point every defined name whose range lived at/after column N so it
#    tracks the shifted data, and register the new "localdb" name.
$wb.Names.Item("macro").RefersTo    = "='#system'!`$O`$2:`$O`$4"
$wb.Names.Item("mail").RefersTo     = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo   = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo      = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo    = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo    = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo      = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo    = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo      = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo     = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("web").RefersTo      = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo       = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo      = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("target").RefersTo   = "='#system'!`$A`$2:`$A`$30"

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
